$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 0.9804946939899589
$ws.Range("C2").Value = 0.02111941343303251
$ws.Range("F2").Value = [double]"1.002544042467475E-05"
$ws.Range("G2").Value = 0.0007469541419202988
$ws.Range("H2").Value = 0.788173738660096
$ws.Range("I2").Value = 0.6985899269312332
$ws.Range("B3").Value = 0.9802013277456616
$ws.Range("C3").Value = 0.02124921321023904
$ws.Range("G3").Value = 0.0005123343433966052
$ws.Range("H3").Value = 0.746236593610306
$ws.Range("I3").Value = 0.5569000901919152
$ws.Range("B4").Value = 0.00905034109549803
$ws.Range("C4").Value = 0.007082435258745194
$ws.Range("F4").Value = [double]"7.495979388161009E-05"
$ws.Range("G4").Value = 0.0007563116016396883
$ws.Range("H4").Value = 0.2795596909893495
$ws.Range("I4").Value = 0.493126246596423

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 0.9505285992633219
$ws.Range("C2").Value = 0.0481628733105622
$ws.Range("G2").Value = 0.001991755136295847
$ws.Range("H2").Value = 0.9935169875817187
$ws.Range("I2").Value = 1.029937367168262
$ws.Range("B3").Value = 0.9510330568296321
$ws.Range("C3").Value = 0.04840547095840617
$ws.Range("G3").Value = 0.001969193055939746
$ws.Range("H3").Value = 0.9873406510393459
$ws.Range("I3").Value = 0.9748426250180053
$ws.Range("B4").Value = 0.009685979398901486
$ws.Range("C4").Value = 0.004094414530504854
$ws.Range("G4").Value = 0.0009343104559375776
$ws.Range("H4").Value = 0.2080728372731967
$ws.Range("I4").Value = 0.4235300205889955

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 0.9210671452578729
$ws.Range("C2").Value = 0.07484834670542642
$ws.Range("D2").Value = 0.0006767272104340685
$ws.Range("F2").Value = 0.0005574211738860795
$ws.Range("G2").Value = 0.002284457111328019
$ws.Range("H2").Value = 0.5900929104719689
$ws.Range("I2").Value = 0.4174617901039237
$ws.Range("B3").Value = 0.9202187954616341
$ws.Range("C3").Value = 0.07484162395487856
$ws.Range("G3").Value = 0.002292001681994532
$ws.Range("H3").Value = 0.5343034131002478
$ws.Range("I3").Value = 0.2855029087167756
$ws.Range("B4").Value = 0.006749509782455811
$ws.Range("C4").Value = 0.003487076506787652
$ws.Range("D4").Value = 0.001264929448071724
$ws.Range("F4").Value = 0.001400371349777551
$ws.Range("G4").Value = 0.001160795321446928
$ws.Range("H4").Value = 0.2644837684096441
$ws.Range("I4").Value = 0.3874999892296863

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 0.8909467915408646
$ws.Range("C2").Value = 0.09826981718768593
$ws.Range("D2").Value = 0.000407237501591936
$ws.Range("F2").Value = 0.0001590283152675463
$ws.Range("G2").Value = 0.007098195311667276
$ws.Range("H2").Value = 0.6144091017493711
$ws.Range("I2").Value = 0.4144231713729987
$ws.Range("B3").Value = 0.8910814951483663
$ws.Range("C3").Value = 0.09796615962870754
$ws.Range("G3").Value = 0.007137929119542013
$ws.Range("H3").Value = 0.588327194415226
$ws.Range("I3").Value = 0.3461302622541393
$ws.Range("B4").Value = 0.005360470018542446
$ws.Range("C4").Value = 0.004782075550153608
$ws.Range("D4").Value = 0.001117154361664795
$ws.Range("F4").Value = 0.0005385346502163977
$ws.Range("G4").Value = 0.0006590618484301129
$ws.Range("H4").Value = 0.1931258736975551
$ws.Range("I4").Value = 0.2564776693969236

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 0.8768849346620484
$ws.Range("C2").Value = 0.1236260336093482
$ws.Range("F2").Value = [double]"3.030826661315419E-05"
$ws.Range("G2").Value = 0.001574482790173337
$ws.Range("H2").Value = 0.6757044307535027
$ws.Range("I2").Value = 0.4967472555568662
$ws.Range("B3").Value = 0.8774334371100203
$ws.Range("C3").Value = 0.1239447215296083
$ws.Range("G3").Value = 0.001482346584668842
$ws.Range("H3").Value = 0.6871177980779859
$ws.Range("I3").Value = 0.4721684669707958
$ws.Range("B4").Value = 0.005839523849284185
$ws.Range("C4").Value = 0.005777157904842808
$ws.Range("F4").Value = 0.0002059005774978092
$ws.Range("G4").Value = 0.0009827128755500966
$ws.Range("H4").Value = 0.2014362014371931
$ws.Range("I4").Value = 0.2755773972803012

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = 0.8086143969765271
$ws.Range("C2").Value = 0.1602747648171917
$ws.Range("D2").Value = 0.001124212936524953
$ws.Range("E2").Value = 0.02722234983353262
$ws.Range("F2").Value = 0.0003254417992280992
$ws.Range("G2").Value = 0.001898146728874546
$ws.Range("H2").Value = 0.5086298223935186
$ws.Range("I2").Value = 0.3030365732268881
$ws.Range("B3").Value = 0.8085324920435435
$ws.Range("C3").Value = 0.1612849298464366
$ws.Range("D3").Value = 0.000395382931832164
$ws.Range("E3").Value = 0.02757040876659583
$ws.Range("G3").Value = 0.001921677985998357
$ws.Range("H3").Value = 0.4857835767625893
$ws.Range("I3").Value = 0.2359877386722556
$ws.Range("B4").Value = 0.01340511838106049
$ws.Range("C4").Value = 0.006807665177567592
$ws.Range("D4").Value = 0.001556251167319811
$ws.Range("E4").Value = 0.01285893253896804
$ws.Range("F4").Value = 0.001184019271675843
$ws.Range("G4").Value = 0.0007453497247259198
$ws.Range("H4").Value = 0.2116130378227949
$ws.Range("I4").Value = 0.233165306699559

$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = 0.7225276718166469
$ws.Range("C2").Value = 0.1683476697025407
$ws.Range("E2").Value = 0.1056918531490265
$ws.Range("F2").Value = 0.0001505186093799428
$ws.Range("G2").Value = 0.001713610567602463
$ws.Range("H2").Value = 0.6032203325446551
$ws.Range("I2").Value = 0.4055468143401456
$ws.Range("B3").Value = 0.7209066823849071
$ws.Range("C3").Value = 0.1679986929826447
$ws.Range("E3").Value = 0.1039545710296024
$ws.Range("G3").Value = 0.001790593826712867
$ws.Range("H3").Value = 0.5774177910919771
$ws.Range("I3").Value = 0.3334140659682546
$ws.Range("B4").Value = 0.01461677138865287
$ws.Range("C4").Value = 0.006153242681700203
$ws.Range("E4").Value = 0.01476287531045645
$ws.Range("F4").Value = 0.0005586456350725643
$ws.Range("G4").Value = 0.0009646668451959455
$ws.Range("H4").Value = 0.2051657244516217
$ws.Range("I4").Value = 0.2562449790239522

$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = 0.708181185005951
$ws.Range("C2").Value = 0.1724371912545334
$ws.Range("E2").Value = 0.1105563555047107
$ws.Range("F2").Value = 0.001345672397094279
$ws.Range("G2").Value = 0.001273726636229345
$ws.Range("H2").Value = 0.5559641398835238
$ws.Range("I2").Value = 0.351177781255972
$ws.Range("B3").Value = 0.7110669465693369
$ws.Range("C3").Value = 0.1728993845819258
$ws.Range("E3").Value = 0.1101695916822096
$ws.Range("F3").Value = 0.0005424366430353156
$ws.Range("G3").Value = 0.00108145737859404
$ws.Range("H3").Value = 0.5432614572229895
$ws.Range("I3").Value = 0.2951336803803157
$ws.Range("B4").Value = 0.01360678609333603
$ws.Range("C4").Value = 0.004975078873867124
$ws.Range("E4").Value = 0.01321033080609342
$ws.Range("F4").Value = 0.001810870680492265
$ws.Range("G4").Value = 0.001092250740575729
$ws.Range("H4").Value = 0.2061715878973304
$ws.Range("I4").Value = 0.2392256851950242

$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = 0.6673998614153978
$ws.Range("C2").Value = 0.1658695130496811
$ws.Range("E2").Value = 0.1632668539891629
$ws.Range("F2").Value = 0.0006022375336771106
$ws.Range("G2").Value = 0.001981173893357759
$ws.Range("H2").Value = 0.3792595324618012
$ws.Range("I2").Value = 0.1621825878764195
$ws.Range("B3").Value = 0.6681029752839038
$ws.Range("C3").Value = 0.1657542500933565
$ws.Range("E3").Value = 0.1622541690094283
$ws.Range("G3").Value = 0.001897682217193756
$ws.Range("H3").Value = 0.381386846505863
$ws.Range("I3").Value = 0.1454641926694007
$ws.Range("B4").Value = 0.01217356199633759
$ws.Range("C4").Value = 0.005193375695966008
$ws.Range("E4").Value = 0.01267891167490279
$ws.Range("F4").Value = 0.001127419942555604
$ws.Range("G4").Value = 0.0005055798376363227
$ws.Range("H4").Value = 0.1361252947544848
$ws.Range("I4").Value = 0.1045827009849272
